$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.579.25"
$ws.Range("E2").Value = "  -1.04%  "
$ws.Range("D3").Value = "1.596.47"
$ws.Range("E3").Value = "  -2.05%  "
$ws.Range("E4").Value = "  +0.56%  "
$ws.Range("D5").Value = "'208.16"
$ws.Range("E5").Value = "  -1.38%  "
$ws.Range("D6").Value = "'0.503"
$ws.Range("E6").Value = "  -3.46%  "
$ws.Range("E7").Value = "  +0.59%  "
$ws.Range("D8").Value = "'22.29"
$ws.Range("E8").Value = "  -4.55%  "
$ws.Range("E9").Value = "  -1.94%  "
$ws.Range("E10").Value = "  -3.28%  "
$ws.Range("D11").Value = "'0.0864"
$ws.Range("E11").Value = "  -1.69%  "
$ws.Range("D12").Value = "1.824.79"
$ws.Range("E12").Value = "  -1.92%  "
$ws.Range("D13").Value = "1.585.65"
$ws.Range("E13").Value = "  -2.71%  "
$ws.Range("E14").Value = "  -3.98%  "
$ws.Range("D15").Value = "'0.536"
$ws.Range("E15").Value = "  -4.54%  "
$ws.Range("D16").Value = "'63.48"
$ws.Range("E16").Value = "  -2.74%  "
$ws.Range("D17").Value = "27.598.00"
$ws.Range("E17").Value = "  -0.99%  "
$ws.Range("D18").Value = "'217.38"
$ws.Range("E18").Value = "  -5.34%  "
$ws.Range("D19").Value = "'7.39"
$ws.Range("E19").Value = "  -3.49%  "
$ws.Range("D20").Value = "0.0₃0694"
$ws.Range("E20").Value = "  -3.65%  "
$ws.Range("E21").Value = "  +0.52%  "
$ws.Range("D22").Value = "'4.19"
$ws.Range("E22").Value = "  -3.33%  "
$ws.Range("D23").Value = "'9.66"
$ws.Range("E23").Value = "  -4.31%  "
$ws.Range("E24").Value = "  -2.52%  "
$ws.Range("D25").Value = "'153.34"
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("E26").Value = "  +0.57%  "
$ws.Range("D27").Value = "'6.73"
$ws.Range("E27").Value = "  -2.33%  "
$ws.Range("E28").Value = "  -2.81%  "
$ws.Range("E29").Value = "  -3.93%  "
$ws.Range("E30").Value = "  -1.31%  "
$ws.Range("D31").Value = "'0.0467"
$ws.Range("E31").Value = "  -2.96%  "
$ws.Range("E32").Value = "  -4.23%  "
$ws.Range("D33").Value = "1.368.58"
$ws.Range("E33").Value = "  -1.81%  "
$ws.Range("D34").Value = "'2.95"
$ws.Range("E34").Value = "  -4.76%  "
$ws.Range("D35").Value = "'1.53"
$ws.Range("E35").Value = "  -3.70%  "
$ws.Range("E36").Value = "  -3.85%  "
$ws.Range("E37").Value = "  -0.84%  "
$ws.Range("E38").Value = "  -2.54%  "
$ws.Range("D39").Value = "'0.541"
$ws.Range("E39").Value = "  -2.84%  "
$ws.Range("D40").Value = "'0.812"
$ws.Range("E40").Value = "  -4.65%  "
$ws.Range("E41").Value = "  +0.55%  "
$ws.Range("E42").Value = "  -3.56%  "
$ws.Range("D43").Value = "'5.37"
$ws.Range("E43").Value = "  -1.01%  "
$ws.Range("D44").Value = "'1.78"
$ws.Range("D45").Value = "'64.04"
$ws.Range("E45").Value = "  -2.73%  "
$ws.Range("D46").Value = "1.734.92"
$ws.Range("E46").Value = "  -1.92%  "
$ws.Range("D47").Value = "'2.12"
$ws.Range("E47").Value = "  -1.64%  "
$ws.Range("D48").Value = "'87.99"
$ws.Range("E48").Value = "  -0.28%  "
$ws.Range("E49").Value = "  -3.03%  "
$ws.Range("D50").Value = "'0.0970"
$ws.Range("E50").Value = "  -4.51%  "
$ws.Range("E51").Value = "  -0.93%  "
